$p = $ppt.ActivePresentation

# Slide 2: Title question text update
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Question: Should ZoomInfo purchase the vendor data?"

# Slide 3: Title question text update (same change as slide 2)
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Question: Should ZoomInfo purchase the vendor data?"

# Slide 3: Content placeholder - update third paragraph text
$contentShape = $s3.Shapes.Item(2)
$contentShape.TextFrame.TextRange.Paragraphs(3).Runs(1).Text = "The comments detail the process"
